# Find the shape whose text contains the phrase we need to change
# ("Call Mr. Taylor / Ms. Mac" -> "Call Mr. Taylor / Mr. K"), instead of
# hard-coding slide/shape indices.
$p = $ppt.ActivePresentation

$oldPhrase = "Call Mr. Taylor / Ms. Mac"
$newPhrase = "Call Mr. Taylor / Mr. K"

$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $shapeRange = $shape.TextFrame.TextRange
            if ($shapeRange.Text -like "*$oldPhrase*") {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the exact paragraph holding the phrase.
$paraCount = $tr.Paragraphs().Count
$targetParaStart = 0
$targetParaLen = 0
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -eq $oldPhrase) {
        $targetParaStart = $para.Start
        $targetParaLen = $para.Length
    }
}

# Rewrite the whole paragraph span in a single assignment so the run isn't
# split anywhere inside it (this keeps the formerly-single run intact with
# its original rPr/dirty state, matching the "Call Mr. Taylor / " prefix
# that is left untouched by the real edit).
$fullSpan = $tr.Characters($targetParaStart, $targetParaLen)
$fullSpan.Text = $newPhrase

# Split off the trailing "K" into its own run by re-assigning that single
# character in place (the paragraph shrank from 25 to $newPhrase.Length
# characters) - this reproduces the two-run split seen in the target edit:
#   run 1: "Call Mr. Taylor / Mr. "
#   run 2: "K"
$newParaLen = $newPhrase.Length
$lastChar = $tr.Characters($targetParaStart + $newParaLen - 1, 1)
$lastChar.Text = "K"

Write-Host "Updated text: $($tr.Text)"
